$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row 11 with the "Nom des images pas adaptés" audit entry,
# reusing the same values/styles pattern as the other rows in the table (row 10).
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)
$ws.Range("F9").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "SEO"
$ws.Range("B11").Value = "Nom des images pas adaptés"
$ws.Range("C11").Value = "Le nom des images doit être explicite et allez droit au but pour que Google puisse comprendre de quoi il s'agit, or ce n'est pas le cas."
$ws.Range("D11").Value = "Mettre des mots décrivant de manière claire ce qu'est l'image."
$ws.Range("E11").Value = "X"

$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.anthedesign.fr/creation-de-sites-internet/format-webp/", [System.Type]::Missing, [System.Type]::Missing, "anthedesign.fr")

# Update the active sheet view to match the new selection/scroll position.
$ws.Range("C1").Select()
$ws.Range("E13").Select()
